$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.416.01"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "1.826.78"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4466"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3759"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07403"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8726"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.828.40"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.706"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.414"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07086"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.0000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008814"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "27.416.41"
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.329"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.953"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.339"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08896"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7915"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.89%  "
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.546"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.946"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9994"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.100"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05264"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.376"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5336"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.347"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.05%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.869"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1701"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.669"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5079"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.683"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9991"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06380"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.20%  "
